$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 468, shifting existing rows 468-555 down to 469-556.
$ws.Rows.Item(468).Insert()

# Populate the newly inserted row 468 with its values.
$ws.Cells.Item(468, 1).Value = 9
$ws.Cells.Item(468, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(468, 3).Value = "Metropolitana"
$ws.Cells.Item(468, 4).Value = 45015
$ws.Cells.Item(468, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(468, 5).Value = 13
$ws.Cells.Item(468, 6).Value = 100112039
$ws.Cells.Item(468, 7).Value = "Ciboulette"
$ws.Cells.Item(468, 8).Value = "Sin especificar"
$ws.Cells.Item(468, 9).Value = "Primera"
$ws.Cells.Item(468, 10).Value = 340
$ws.Cells.Item(468, 11).Value = 1000
$ws.Cells.Item(468, 12).Value = 1200
$ws.Cells.Item(468, 13).Value = 1100
$ws.Cells.Item(468, 14).Value = "$/docena de atados"
$ws.Cells.Item(468, 15).Value = "Región Metropolitana"
$ws.Cells.Item(468, 16).Value = 367
$ws.Cells.Item(468, 17).Value = 3
$ws.Cells.Item(468, 18).Value = "Hortaliza"
